$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text so numeric-looking
# strings (e.g. "69.379.84", "0.0901", "6.40") keep their exact
# textual representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.379.84"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "3.679.16"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "685.37"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("D6").Value = "158.97"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -1.05%  "

# Row 9
$ws.Range("E9").Value = "  -1.54%  "

# Row 10
$ws.Range("D10").Value = "7.04"
$ws.Range("E10").Value = "  -3.43%  "

# Row 11
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -3.20%  "

# Row 12
$ws.Range("E12").Value = "  -1.22%  "

# Row 13
$ws.Range("D13").Value = "4.298.06"
$ws.Range("E13").Value = "  -0.34%  "

# Row 14
$ws.Range("D14").Value = "32.32"
$ws.Range("E14").Value = "  -3.70%  "

# Row 15
$ws.Range("D15").Value = "3.672.57"
$ws.Range("E15").Value = "  -0.47%  "

# Row 16
$ws.Range("D16").Value = "69.336.19"
$ws.Range("E16").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "15.85"
$ws.Range("E18").Value = "  -3.01%  "

# Row 19
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -3.23%  "

# Row 20
$ws.Range("D20").Value = "469.80"
$ws.Range("E20").Value = "  -2.37%  "

# Row 21
$ws.Range("D21").Value = "9.98"
$ws.Range("E21").Value = "  +1.78%  "

# Row 22
$ws.Range("E22").Value = "  -2.52%  "

# Row 23
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("D24").Value = "3.823.98"
$ws.Range("E24").Value = "  -0.29%  "

# Row 25
$ws.Range("E25").Value = "  +0.15%  "

# Row 26
$ws.Range("E26").Value = "  -4.47%  "

# Row 27
$ws.Range("D27").Value = "10.96"
$ws.Range("E27").Value = "  -4.78%  "

# Row 28
$ws.Range("D28").Value = "9.17"
$ws.Range("E28").Value = "  -4.37%  "

# Row 29
$ws.Range("E29").Value = "  -1.28%  "

# Row 30
$ws.Range("E30").Value = "  -5.61%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "6.55"
$ws.Range("E31").Value = "  -4.21%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -5.99%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "26.90"
$ws.Range("E33").Value = "  -0.55%  "

# Row 34
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.34%  "

# Row 35
$ws.Range("D35").Value = "3.653.30"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("E36").Value = "  -4.27%  "

# Row 37
$ws.Range("D37").Value = "8.14"
$ws.Range("E37").Value = "  -4.63%  "

# Row 38
$ws.Range("E38").Value = "  +2.05%  "

# Row 39
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("D40").Value = "2.22"
$ws.Range("E40").Value = "  +1.47%  "

# Row 41
$ws.Range("D41").Value = "0.0901"
$ws.Range("E41").Value = "  -4.36%  "

# Row 42
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("E43").Value = "  -2.30%  "

# Row 44
$ws.Range("D44").Value = "166.21"
$ws.Range("E44").Value = "  +4.82%  "

# Row 45
$ws.Range("E45").Value = "  -1.49%  "

# Row 46
$ws.Range("D46").Value = "0.000283"
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -4.83%  "

# Row 48
$ws.Range("E48").Value = "  +4.64%  "

# Row 49
$ws.Range("E49").Value = "  +0.10%  "

# Row 50
$ws.Range("D50").Value = "7.79"
$ws.Range("E50").Value = "  -3.83%  "

# Row 51
$ws.Range("D51").Value = "27.22"
$ws.Range("E51").Value = "  -1.20%  "
